# Updated cryptos list on Fri Jun 21 15:31:19 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = "new price text"; E = "new volume text" }
$updates = @{
    2  = @{ D = "63.911.99";  E = "  -1.45%  " }
    3  = @{ D = "3.498.40";   E = "  -0.70%  " }
    5  = @{ D = "583.42";     E = "  -2.19%  " }
    6  = @{ D = "132.00";     E = "  -1.61%  " }
    7  = @{ D = "3.496.79";   E = "  -0.72%  " }
    8  = @{ E = "  -0.02%  " }
    9  = @{ D = "0.486";      E = "  -1.94%  " }
    10 = @{ D = "0.123";      E = "  -0.34%  " }
    11 = @{ E = "  +0.33%  " }
    12 = @{ D = "0.381";      E = "  -0.38%  " }
    13 = @{ D = "4.097.25";   E = "  -0.64%  " }
    14 = @{ D = "27.58";      E = "  +1.14%  " }
    15 = @{ E = "  +1.75%  " }
    16 = @{ D = "0.0000178";  E = "  -2.05%  " }
    17 = @{ D = "3.479.83";   E = "  -1.23%  " }
    18 = @{ D = "64.059.67";  E = "  -1.32%  " }
    19 = @{ D = "9.99";       E = "  +0.32%  " }
    20 = @{ D = "14.35";      E = "  -0.37%  " }
    21 = @{ D = "5.68";       E = "  +0.01%  " }
    22 = @{ D = "386.31";     E = "  -1.37%  " }
    23 = @{ D = "0.579";      E = "  +0.57%  " }
    24 = @{ D = "3.641.85";   E = "  -0.68%  " }
    25 = @{ D = "73.05";      E = "  -1.61%  " }
    27 = @{ D = "0.0000113";  E = "  +0.12%  " }
    28 = @{ E = "  -1.92%  " }
    29 = @{ D = "7.54";       E = "  -2.56%  " }
    30 = @{ D = "1.00";       E = "  -0.20%  " }
    31 = @{ D = "2.25";       E = "  -1.40%  " }
    32 = @{ D = "8.29" }
    33 = @{ D = "3.506.27";   E = "  -0.61%  " }
    34 = @{ E = "  -0.02%  " }
    35 = @{ D = "23.61";      E = "  -1.99%  " }
    36 = @{ E = "  +0.09%  " }
    37 = @{ D = "5.37";       E = "  +2.67%  " }
    38 = @{ D = "1.58";       E = "  +0.44%  " }
    39 = @{ D = "6.96";       E = "  +1.99%  " }
    40 = @{ D = "162.77";     E = "  -3.77%  " }
    41 = @{ D = "0.0799";     E = "  -2.64%  " }
    42 = @{ D = "26.52";      E = "  +5.10%  " }
    43 = @{ D = "0.811";      E = "  -1.25%  " }
    44 = @{ E = "  +0.03%  " }
    45 = @{ E = "  -0.94%  " }
    46 = @{ D = "41.59";      E = "  -2.45%  " }
    47 = @{ D = "4.40";       E = "  -0.44%  " }
    48 = @{ D = "1.64";       E = "  -0.50%  " }
    49 = @{ D = "6.86";       E = "  -0.72%  " }
    50 = @{ D = "2.421.92";   E = "  +0.78%  " }
    51 = @{ D = "0.896";      E = "  +0.00%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals["D"]
        $cell.Style = "Normal"
    }
    if ($vals.ContainsKey("E")) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals["E"]
        $cell.Style = "Normal"
    }
}
